# Add a new "2021" column (M) to the mobile-network-coverage table, mirroring
# the formatting of the existing 2020 column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the 2020 column (L, rows 3-7) onto the new column
# (M, rows 3-7) so the new cells match the surrounding table style.
$ws.Range("L3:L7").Copy()
$ws.Range("M3:M7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header: year 2021
$ws.Range("M4").Value = 2021

# Data rows (2G / 3G / 4G coverage percentages for 2021)
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96

# Reset the view back to the top-left cell (the sheet had been left scrolled
# to B1 with N13 selected).
$ws.Range("A1").Select()
